$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.5
$ws.Range("I3").Value = 5.5
$ws.Range("J3").Value = 2
$ws.Range("L3").Value = 5.5
$ws.Range("Q3").Value = 1.53
$ws.Range("R3").Value = 2.4
$ws.Range("Z3").Value = 11
$ws.Range("AA3").Value = 11
$ws.Range("AD3").Value = 9
$ws.Range("AM3").Value = 41
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 7.5
$ws.Range("AX3").Value = 29
$ws.Range("BA3").Value = 101
